$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append today's time log as row 8: date + start/break-out/break-in/end times
# (10:00, 12:30, 13:00, 13:30), matching the pattern of the existing rows.
$ws.Range("A8").Value = "16.9.2025"
$ws.Range("B8").Value = 10 / 24
$ws.Range("C8").Value = 12.5 / 24
$ws.Range("D8").Value = 13 / 24
$ws.Range("E8").Value = 13.5 / 24

# Match the time number format used by the existing rows (numFmtId 18, style index 1).
$ws.Range("B8:E8").NumberFormat = $ws.Range("B2:E2").NumberFormat

# Move the selection the way Excel would after entering data through row 8.
$ws.Range("F10").Select()
